$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q4" right before the existing
#    "2022-Q3" sheet, so the tab order becomes:
#      总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with the fund-holding detail rows,
#    matching the layout used by its sibling quarter sheets.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

$q4rows = @(
    @(0, "090019", "大成景恒混合A", "1.18", "93.72", "1.63", "0.0192", 7),
    @(1, "006038", "大成景恒混合C", "0.89", "93.72", "1.63", "0.0145", 7)
)

for ($r = 0; $r -lt $q4rows.Length; $r++) {
    $row = $q4rows[$r]
    $excelRow = $r + 2

    $q4.Cells.Item($excelRow, 1).Value = $row[0]
    $q4.Cells.Item($excelRow, 2).Value = $row[1]
    $q4.Cells.Item($excelRow, 3).Value = $row[2]
    $q4.Cells.Item($excelRow, 4).Value = $row[3]
    $q4.Cells.Item($excelRow, 5).Value = $row[4]
    $q4.Cells.Item($excelRow, 6).Value = $row[5]
    $q4.Cells.Item($excelRow, 7).Value = $row[6]
    $q4.Cells.Item($excelRow, 8).Value = $row[7]

    $aCell = $q4.Cells.Item($excelRow, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift the existing quarters down one
#    row and insert the new 2022-Q4 totals at the top of the data.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2022-Q1"
$total.Cells.Item(5, 3).Value = 4
$total.Cells.Item(5, 4).Value = 1.05

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 4
$total.Cells.Item(4, 4).Value = 1.3

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 6
$total.Cells.Item(3, 4).Value = 0.66

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.03

# Row 5 column A is brand new territory on this sheet - copy the formatting
# (bold/border/alignment) that the other "id" cells in column A already use.
$total.Cells.Item(4, 1).Copy()
$total.Cells.Item(5, 1).PasteSpecial(-4122)   # xlPasteFormats
